$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Hunk 1: merge the two "TUE Apr 16" / " 14:16:46 IST 2019" runs into a
# single run containing "TUE Apr 16 14:16:46 IST 2019".
# ---------------------------------------------------------------------
$d.Content.Find.Execute("TUE Apr 16 14:16:46 IST 2019", $false, $false, $false, $false, $false, $true, 1, $false, "TUE Apr 16 14:16:46 IST 2019", 2) | Out-Null

# ---------------------------------------------------------------------
# Hunk 2: append a brand new "purchase" record right after the final
# "Amount Received mode ... - CASH" paragraph (the one closing the
# previous block, right before the trailing blank paragraphs at the
# end of the document).
# ---------------------------------------------------------------------

# Locate the LAST paragraph whose text is "Amount Received mode" + tabs
# + "- CASH" (not "- CASH AND CLEARD" / "- CASH AND CLEAR").
$matches = @()
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("Amount Received mode") -and $t.EndsWith("- CASH`r")) {
        $matches += $p
    }
}
$cur = $matches[$matches.Count - 1]

# --- New paragraph 1: empty, bold (Courier New) -----------------------
$r = $cur.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$cur = $cur.Next()

# --- New paragraph 2: "SUN Apr 21" + " 14:29:33 IST 2019" (not bold) --
$r = $cur.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$cur = $cur.Next()
$r = $cur.Range
$r.Font.Bold = $false
$r.InsertAfter("SUN Apr 21")
$cur2 = $cur.Range
$cur2.Collapse(0)
$cur2.Font.Bold = $false
$cur2.InsertAfter(" 14:29:33 IST 2019")

# --- New paragraph 3: "Person Name" ... "- CS" -------------------------
$r = $cur.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$cur = $cur.Next()
$r = $cur.Range
$r.Font.Bold = $false
$r.InsertAfter("Person Name`t`t`t`t- CS")

# --- New paragraph 4: "Bill number" ... "- 12344" -----------------------
$r = $cur.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$cur = $cur.Next()
$r = $cur.Range
$r.Font.Bold = $false
$r.InsertAfter("Bill number`t`t`t`t- 12344")

# --- New paragraph 5: separator line -----------------------------------
$r = $cur.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$cur = $cur.Next()
$r = $cur.Range
$r.Font.Bold = $false
$r.InsertAfter("---------------------------------------------------------------")

# --- New paragraph 6: "Item Name" ... "- CARROT" ------------------------
$r = $cur.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$cur = $cur.Next()
$r = $cur.Range
$r.Font.Bold = $false
$r.InsertAfter("Item Name`t`t`t`t- CARROT")

# --- New paragraph 7: "Number of Pockets" ... "- 1" ----------------------
$r = $cur.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$cur = $cur.Next()
$r = $cur.Range
$r.Font.Bold = $false
$r.InsertAfter("Number of Pockets`t`t`t- 1")

# --- New paragraph 8: "Number of KGs" ... "- 94" --------------------------
$r = $cur.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$cur = $cur.Next()
$r = $cur.Range
$r.Font.Bold = $false
$r.InsertAfter("Number of KGs`t`t`t- 94")

# --- New paragraph 9: "Rate" ... "- 24" -----------------------------------
$r = $cur.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$cur = $cur.Next()
$r = $cur.Range
$r.Font.Bold = $false
$r.InsertAfter("Rate`t`t`t`t`t- 24")

# --- New paragraph 10: "Total Price" ... "- 2256.0" -----------------------
$r = $cur.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$cur = $cur.Next()
$r = $cur.Range
$r.Font.Bold = $false
$r.InsertAfter("Total Price`t`t`t`t- 2256.0")

# --- New paragraph 11: "Amount balance" ... "- 3356.0" (bold) -------------
$r = $cur.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$cur = $cur.Next()
$r = $cur.Range
$r.Font.Bold = $true
$r.InsertAfter("Amount balance`t`t`t- 3356.0")

# --- New paragraph 12: empty, bold ----------------------------------------
$r = $cur.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Font.Bold = $true

# --- New paragraph 13: empty, bold (shares closing tag with the pre------
#     existing first trailing blank paragraph in the original document) --
$r = $cur.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.Font.Bold = $true
